$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '49.473.82'
$ws.Range('E2').Value = '  -0.89%  '

# Row 3
$ws.Range('D3').Value = '2.632.06'
$ws.Range('E3').Value = '  -0.86%  '

# Row 4
$ws.Range('E4').Value = '  +0.07%  '

# Row 5
$ws.Range('D5').Value = '''111.30'
$ws.Range('E5').Value = '  -2.02%  '

# Row 6
$ws.Range('D6').Value = '''325.42'
$ws.Range('E6').Value = '  -0.81%  '

# Row 8
$ws.Range('D8').Value = '''1.00'
$ws.Range('E8').Value = '  +0.04%  '

# Row 9
$ws.Range('E9').Value = '  -1.79%  '

# Row 10
$ws.Range('D10').Value = '''39.41'
$ws.Range('E10').Value = '  -4.44%  '

# Row 11
$ws.Range('D11').Value = '''20.15'
$ws.Range('E11').Value = '  -0.11%  '

# Row 12
$ws.Range('E12').Value = '  -1.79%  '

# Row 13
$ws.Range('E13').Value = '  +1.43%  '

# Row 14
$ws.Range('D14').Value = '''7.37'
$ws.Range('E14').Value = '  +0.16%  '

# Row 15
$ws.Range('D15').Value = '3.044.34'
$ws.Range('E15').Value = '  -0.77%  '

# Row 16
$ws.Range('D16').Value = '2.633.06'
$ws.Range('E16').Value = '  -1.35%  '

# Row 17
$ws.Range('E17').Value = '  -2.64%  '

# Row 18
$ws.Range('D18').Value = '49.440.78'
$ws.Range('E18').Value = '  -0.82%  '

# Row 19
$ws.Range('D19').Value = '''13.05'
$ws.Range('E19').Value = '  -0.91%  '

# Row 20
$ws.Range('B20').Value = 'ImmutableX'
$ws.Range('C20').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D20').Value = '''2.90'
$ws.Range('E20').Value = '  -1.32%  '

# Row 21
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').Value = '''6.67'
$ws.Range('E21').Value = '  -1.55%  '

# Row 22
$ws.Range('E22').Value = '  -1.31%  '

# Row 23
$ws.Range('D23').Value = '''267.81'
$ws.Range('E23').Value = '  -3.26%  '

# Row 24
$ws.Range('D24').Value = '''68.89'
$ws.Range('E24').Value = '  -4.70%  '

# Row 25
$ws.Range('E25').Value = '  -2.15%  '

# Row 26
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').Value = '''1.00'
$ws.Range('E26').Value = '  +0.15%  '

# Row 27
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '''25.99'
$ws.Range('E27').Value = '  -3.34%  '

# Row 28
$ws.Range('D28').Value = '''10.14'
$ws.Range('E28').Value = '  +1.66%  '

# Row 29
$ws.Range('E29').Value = '  -1.14%  '

# Row 30
$ws.Range('E30').Value = '  -1.50%  '

# Row 31
$ws.Range('D31').Value = '''34.54'
$ws.Range('E31').Value = '  -4.03%  '

# Row 32
$ws.Range('D32').Value = '''49.54'
$ws.Range('E32').Value = '  -1.41%  '

# Row 33
$ws.Range('D33').Value = '''5.46'
$ws.Range('E33').Value = '  +0.78%  '

# Row 34
$ws.Range('D34').Value = '''0.0808'
$ws.Range('E34').Value = '  +0.59%  '

# Row 35
$ws.Range('E35').Value = '  -0.08%  '

# Row 36
$ws.Range('D36').Value = '''18.98'
$ws.Range('E36').Value = '  -3.04%  '

# Row 37
$ws.Range('D37').Value = '''4.95'
$ws.Range('E37').Value = '  +3.49%  '

# Row 39
$ws.Range('E39').Value = '  -0.01%  '

# Row 40
$ws.Range('D40').Value = '''128.93'
$ws.Range('E40').Value = '  +2.63%  '

# Row 41
$ws.Range('D41').Value = '''22.90'
$ws.Range('E41').Value = '  +2.77%  '

# Row 42
$ws.Range('E42').Value = '  -1.61%  '

# Row 43
$ws.Range('D43').Value = '''2.20'
$ws.Range('E43').Value = '  -1.27%  '

# Row 44
$ws.Range('D44').Value = '''0.0326'
$ws.Range('E44').Value = '  +3.36%  '

# Row 45
$ws.Range('D45').Value = '2.041.65'
$ws.Range('E45').Value = '  -1.38%  '

# Row 46
$ws.Range('D46').Value = '''2.16'
$ws.Range('E46').Value = '  +8.70%  '

# Row 47
$ws.Range('E47').Value = '  -3.69%  '

# Row 48
$ws.Range('E48').Value = '  -4.54%  '

# Row 49
$ws.Range('D49').Value = '''8.83'
$ws.Range('E49').Value = '  -3.59%  '

# Row 50
$ws.Range('E50').Value = '  -3.72%  '

# Row 51
$ws.Range('D51').Value = '''58.45'
$ws.Range('E51').Value = '  +1.00%  '
